$d = $word.ActiveDocument

$replacements = @(
    @("🎯 Categorization Framework", "**Goal:** Categorization Framework"),
    @("📊 Component Classification Matrix", "**Summary:** Component Classification Matrix"),
    @("**FREE Components** ✅", "**FREE Components** **Success:**"),
    @("**LOCKED Components** 🔒", "**LOCKED Components** **Security:**"),
    @("🏗️ Hybrid Components (Context-Aware)", "**Architecture:** Hybrid Components (Context-Aware)"),
    @("**Button Group** 🔄", "**Button Group** **Process:**"),
    @("**Form Fields** 🔄", "**Form Fields** **Process:**"),
    @("**Modal/Dialog** 🔄", "**Modal/Dialog** **Process:**"),
    @("🧩 Nesting Scenarios Analysis", "**Note:** Nesting Scenarios Analysis"),
    @("🎯 Size Inheritance Rules", "**Goal:** Size Inheritance Rules"),
    @("💡 Token Architecture Implications", "**Tip:** Token Architecture Implications"),
    @("🔍 Competitive Validation", "**Analysis:** Competitive Validation")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Update footer date
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers(1)
    $ftr.Range.Find.Execute("2025-09-07 12:53", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-12 17:37", 2)
}
